{"js": "// Locate the last paragraph in the document body (\"Las siglas DAO ...\")\n// and append the two new paragraphs described by the diff right after it,\n// before the final section break.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst firstNew = lastParagraph.insertParagraph(\n  \"Esta clase cuenta con atributos y m\u00e9todos de clase. Los atributos son SELECCIONAR; ELIMINAR; ACTUALIZAR E INSERTAR.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nfirstNew.insertParagraph(\n  \"Los m\u00e9todos, de clase, utilizan estas sentencias, junto con datos tra\u00eddos por argumentos (objetos de clase Persona) para realizar las operaciones CRUD en la DB.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Append the two new paragraphs described by the diff right after the\n# document's last paragraph (\"Las siglas DAO ...\"), before the section break.\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$firstNew = $d.Paragraphs.Last\n$firstNew.Range.Text = \"Esta clase cuenta con atributos y m\u00e9todos de clase. Los atributos son SELECCIONAR; ELIMINAR; ACTUALIZAR E INSERTAR.\"\n\n$firstNew = $d.Paragraphs.Last\n$firstNew.Range.InsertParagraphAfter()\n\n$secondNew = $d.Paragraphs.Last\n$secondNew.Range.Text = \"Los m\u00e9todos, de clase, utilizan estas sentencias, junto con datos tra\u00eddos por argumentos (objetos de clase Persona) para realizar las operaciones CRUD en la DB.\"\n"}
